$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("_itemType"), shifting existing
# columns B:S to C:T.
$ws.Columns("B").Insert()

# New column inherits column A's width (16).
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Header for the newly inserted column.
$ws.Range("B1").Value = "_itemType"

# Fill the new column with the item type for every data row.
$ws.Range("B2:B5").Value = "Equipments"

$ws.Range("F11").Select()
